$wb = $excel.ActiveWorkbook

# Offense (OFF) sheet - row 3 ("R") target depth data update
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 170
$wsOff.Range("C3").Value = 111
$wsOff.Range("D3").Value = 33
$wsOff.Range("E3").Value = 21

# Defense (DEF) sheet - row 3 ("R") target depth data update
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 188
$wsDef.Range("C3").Value = 126
$wsDef.Range("D3").Value = 40
$wsDef.Range("E3").Value = 15
